# Adds new columns I ("I0") and J ("IF") to the stats sheet, mirroring the
# formatting of the existing header/data columns and filling in the values
# that were added in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, border, centered alignment) from the
# existing "IP" header cell (H1) onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats

# --- Data rows (2-39) ------------------------------------------------------
$data = @(
    @(2, 1, 5),
    @(3, 9, 9),
    @(4, 1, 7),
    @(5, 1, 5),
    @(6, 1, 6),
    @(7, 1, 5),
    @(8, 1, 7),
    @(9, 1, 8),
    @(10, 1, 8),
    @(11, 1, 6),
    @(12, 1, 8),
    @(13, 1, 5),
    @(14, 1, 6),
    @(15, 1, 4),
    @(16, 1, 7),
    @(17, 1, 6),
    @(18, 1, 7),
    @(19, 1, 4),
    @(20, 1, 6),
    @(21, 1, 6),
    @(22, 1, 5),
    @(23, 1, 1),
    @(24, 1, 8),
    @(25, 1, 6),
    @(26, 1, 7),
    @(27, 1, 5),
    @(28, 1, 5),
    @(29, 1, 8),
    @(30, 1, 9),
    @(31, 1, 5),
    @(32, 1, 6),
    @(33, 1, 7),
    @(34, 1, 6),
    @(35, 1, 6),
    @(36, 1, 5),
    @(37, 1, 4),
    @(38, 1, 3),
    @(39, 1, 2)
)

foreach ($item in $data) {
    $row = $item[0]
    $i0 = $item[1]
    $iF = $item[2]
    $ws.Cells.Item($row, 9).Value = $i0
    $ws.Cells.Item($row, 10).Value = $iF
}

# Refresh the sheet's used-range dimension to include the new columns.
$ws.Range("A1:J39").Select()
